# issue #5: property land done
#
# Normalizes the OCR-mangled shared-string text (stray internal spaces and
# thousands-separator commas) on every property sheet, and appends the
# standard metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the "land" (土地)
# sheet so it matches the schema already used elsewhere.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 土地 (land)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header row - rename the first column from the Chinese label to the
# canonical "name" field; the rest of the header row is unchanged.
$ws1.Range("B1").Value = "name"

# Clean up the OCR-mangled text in the data row.
$ws1.Range("B2").Value = "新北市永和區永福段08230000地號"
$ws1.Range("F2").Value = "80年06月13日"
$ws1.Range("G2").Value = "第一次登記"

# New metadata columns I:O, values first (K2 looks like a date, so force
# text entry for it), then copy the look of the existing header/data cells
# onto the new ones.
$ws1.Range("I1").Value = "property_category"
$ws1.Range("J1").Value = "category"
$ws1.Range("K1").Value = "date"
$ws1.Range("L1").Value = "legislator_name"
$ws1.Range("M1").Value = "legislator_id"
$ws1.Range("N1").Value = "source_file"
$ws1.Range("O1").Value = "index"

$ws1.Range("K2").NumberFormat = "@"
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
$ws1.Range("K2").Value = "2012-03-01"
$ws1.Range("L2").Value = "林淑芬"
$ws1.Range("M2").Value = 1337
$ws1.Range("N2").Value = "tmp3f851"
$ws1.Range("O2").Value = 15

$ws1.Range("H1").Copy()
$ws1.Range("I1:O1").PasteSpecial(-4122)
$ws1.Range("H2").Copy()
$ws1.Range("I2:O2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet 2: 建物 (building)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "新北市永和區永福段02017000建號"
$ws2.Range("F2").Value = "80年06月13日"
$ws2.Range("G2").Value = "第一次登.記"

# ---------------------------------------------------------------------
# Sheet 3: 汽車 (car)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "toyotarav4rod"
$ws3.Range("E2").Value = "100年03月01曰"

# ---------------------------------------------------------------------
# Sheet 5: 存款 (deposit)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B4").Value = "遠東國際商業銀行永和分行"
$ws5.Range("B5").Value = "遠東國際商業銀行永和分行"

# F4 keeps its original text type ("1，527，000" -> "1527000", still text,
# not a number), so force text entry and then restore the plain look of
# the other amount cells in this column.
$ws5.Range("F4").NumberFormat = "@"
$ws5.Range("F4").Value = "1527000"
$ws5.Range("F5").Copy()
$ws5.Range("F4").PasteSpecial(-4122)
